$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three advice texts (LoiKhuyen) ---
$ws.Range("B2").Value = 'Bạn là người giàu nghị lực và có thể gặp nhiều khó khăn trong cuộc sống. Tuy nhiên "lửa thử vàng, gian nan thử sức", hãy cố gắng trở thành hòn ngọc quý.'
$ws.Range("B3").Value = "Bạn gặp nhiều may mắn hơn thực lực bản thân. Hãy cố gắng trau dồi năng lực của bản thân để càng gặp nhiều điều toại nguyện hơn."
$ws.Range("B8").Value = "Bạn sinh ra gặp rất nhiều sóng gió cuộc đời nhưng đừng nản chí và bỏ cuộc vì đến khi vào đại vận trung niên bạn sẽ được hưởng trọn vẹn thành quả của sự cố gắng do quá trình tích lũy trước đó, cuộc sống gắn liền phần nhiều đến tín ngưỡng và tôn giáo. Bạn hãy cố gắng trau dồi kinh nghiệm, kiến thức và trải nghiệm."

# --- Highlight the section-header cells (A1, A2, A3, A8) with an accent (theme) fill color ---
$ws.Range("A1").Interior.ThemeColor = 10
$ws.Range("A2").Interior.ThemeColor = 10
$ws.Range("A3").Interior.ThemeColor = 10
$ws.Range("A8").Interior.ThemeColor = 10

# --- Move the active selection to A8 ---
[void]$ws.Range("A8").Select()
